$wb = $excel.ActiveWorkbook

# "To fix" is the first sheet (sheetId=1, r:id=rId1) -> sheet1.xml
$ws = $wb.Worksheets.Item("To fix")

# Append the three new bug entries to the bottom of the list (rows 8-10)
$ws.Range("A8").Value = "reloading the page with an loaded exercise should reload the list of assignments page."
$ws.Range("A9").Value = "disable check answer when nothing is selected yet."
$ws.Range("A10").Value = "the right answer is always first. Needs to be fixed."

# Make "To fix" the active/selected sheet and move the selection to the new last cell
$ws.Activate()
$ws.Range("A10").Select()

$wb.Save()
